$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.123.80"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.864.36"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'313.14"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "'0.5094"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("D8").Value = "'0.3896"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "'0.08205"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("B11").Value = "Polkadot"
$ws.Range("C11").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D11").Value = "'6.180"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.863.16"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.16"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.190"
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001094"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'90.56"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06671"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.973"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "28.163.69"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'11.05"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.216"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'3.394"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "2.062.76"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'159.10"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'20.60"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'2.409"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'125.47"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'0.1051"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'1.035"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "'5.834"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'9.295"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02416"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06518"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.2182"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").Value = "'1.234"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'1.173"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "'4.943"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "'11.11"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "'12.93"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").Value = "'3.673"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "'1.271"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'1.986"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").Value = "'1.201"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "'121.21"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'0.06865"
$ws.Range("E51").Value = "  +0.50%  "
